$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 4407.0264
$ws.Range("J17").Value = 4407.0264
$ws.Range("L17").Value = 13221.0792
$ws.Range("N17").Value = -13557.0792
# row 18
$ws.Range("H18").Value = 625
$ws.Range("I18").Value = 625
$ws.Range("K18").Value = 625
$ws.Range("M18").Value = -341
# row 19
$ws.Range("H19").Value = 31250946
$ws.Range("I19").Value = 766
$ws.Range("J19").Value = 62501124
$ws.Range("K19").Value = 766
$ws.Range("L19").Value = 62501124
$ws.Range("M19").Value = -591
$ws.Range("N19").Value = -62501474
# row 33
$ws.Range("H33").Value = 377.3913
$ws.Range("I33").Value = 375
$ws.Range("K33").Value = 375
$ws.Range("M33").Value = -146
# row 137
$ws.Range("H137").Value = 16395942
$ws.Range("I137").Value = 40001980
$ws.Range("K137").Value = 120005940
$ws.Range("M137").Value = -120003390

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 2043401.4
$ws.Range("I32").Value = 2660.2354
$ws.Range("J32").Value = 71428600
$ws.Range("K32").Value = 2660.2354
$ws.Range("L32").Value = 71428600
$ws.Range("M32").Value = -2373.2354
$ws.Range("N32").Value = -71429174
# row 45
$ws.Range("H45").Value = 5534.9
$ws.Range("I45").Value = 2749.8
$ws.Range("K45").Value = 2749.8
$ws.Range("M45").Value = -2372.8
# row 74
$ws.Range("H74").Value = 7410111
$ws.Range("I74").Value = 9011060
$ws.Range("K74").Value = 9011060
$ws.Range("M74").Value = -9010186
# row 77
$ws.Range("H77").Value = 7410111
$ws.Range("I77").Value = 9011060
$ws.Range("K77").Value = 45055300
$ws.Range("M77").Value = -45050932
# row 132
$ws.Range("H132").Value = 2725.6309
$ws.Range("I132").Value = 2092.672
$ws.Range("K132").Value = 6278.016
$ws.Range("M132").Value = -3748.016

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 5
$ws.Range("H5").Value = 1874.4286
$ws.Range("I5").Value = 350.5
$ws.Range("J5").Value = 3906.3333
$ws.Range("K5").Value = 350.5
$ws.Range("L5").Value = 3906.3333
$ws.Range("M5").Value = -237.5
$ws.Range("N5").Value = -4132.3333
# row 7
$ws.Range("H7").Value = 7250.5
$ws.Range("I7").Value = 7250.5
$ws.Range("K7").Value = 7250.5
$ws.Range("M7").Value = -7137.5
# row 22
$ws.Range("H22").Value = 580.75
$ws.Range("I22").Value = 606.5714
$ws.Range("K22").Value = 606.5714
$ws.Range("M22").Value = -433.5714
# row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 2
$ws.Range("H2").Value = 12388.4
$ws.Range("I2").Value = 11027.333
$ws.Range("J2").Value = 14430
$ws.Range("K2").Value = 11027.333
$ws.Range("L2").Value = 14430
$ws.Range("M2").Value = -10914.333
$ws.Range("N2").Value = -14656
# row 132
$ws.Range("H132").Value = 2230.6274
$ws.Range("I132").Value = 1912.6744
$ws.Range("K132").Value = 5738.023200000001
$ws.Range("M132").Value = -3208.023200000001
# row 134
$ws.Range("H134").Value = 4090.7297
$ws.Range("I134").Value = 3595.6292
$ws.Range("K134").Value = 10786.8876
$ws.Range("M134").Value = -8251.8876

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 43
$ws.Range("H43").Value = 4326.6665
$ws.Range("J43").Value = 5490
$ws.Range("L43").Value = 16470
$ws.Range("N43").Value = -16698
# row 81
$ws.Range("H81").Value = 3679.125
$ws.Range("I81").Value = 2772.3572
$ws.Range("J81").Value = 4384.3887
$ws.Range("K81").Value = 8317.071599999999
$ws.Range("L81").Value = 13153.1661
$ws.Range("M81").Value = -7194.071599999999
$ws.Range("N81").Value = -15399.1661
# row 82
$ws.Range("H82").Value = 3061.5
$ws.Range("I82").Value = 3061.5
$ws.Range("K82").Value = 9184.5
$ws.Range("M82").Value = -8778.5
# row 84
$ws.Range("H84").Value = 3679.125
$ws.Range("I84").Value = 2772.3572
$ws.Range("J84").Value = 4384.3887
$ws.Range("K84").Value = 24951.2148
$ws.Range("L84").Value = 39459.49830000001
$ws.Range("M84").Value = -19335.2148
$ws.Range("N84").Value = -50691.49830000001
# row 85
$ws.Range("H85").Value = 3061.5
$ws.Range("I85").Value = 3061.5
$ws.Range("K85").Value = 9184.5
$ws.Range("M85").Value = -7780.5
# row 131
$ws.Range("H131").Value = 5606191
$ws.Range("I131").Value = 13891063
$ws.Range("K131").Value = 41673189
$ws.Range("M131").Value = -41668149
# row 140
$ws.Range("H140").Value = 1912.16
$ws.Range("I140").Value = 1462.4286
$ws.Range("J140").Value = 4273.25
$ws.Range("K140").Value = 4387.2858
$ws.Range("L140").Value = 12819.75
$ws.Range("M140").Value = 792.7142000000003
$ws.Range("N140").Value = -23179.75

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 666.1875
$ws.Range("I2").Value = 69.166664
$ws.Range("J2").Value = 1024.4
$ws.Range("K2").Value = 69.166664
$ws.Range("L2").Value = 1024.4
$ws.Range("M2").Value = 43.833336
$ws.Range("N2").Value = -1250.4
# row 3
$ws.Range("H3").Value = 101719
$ws.Range("I3").Value = 200588.2
$ws.Range("J3").Value = 2849.8
$ws.Range("K3").Value = 200588.2
$ws.Range("L3").Value = 2849.8
$ws.Range("M3").Value = -200472.2
$ws.Range("N3").Value = -3081.8
# row 102
$ws.Range("H102").Value = 3032623.5
$ws.Range("I102").Value = 3335485.8
$ws.Range("K102").Value = 3335485.8
$ws.Range("M102").Value = -3333863.8
# row 107
$ws.Range("H107").Value = 903.375
$ws.Range("J107").Value = 1995.5
$ws.Range("L107").Value = 1995.5
$ws.Range("N107").Value = -5835.5
# row 126
$ws.Range("H126").Value = 5461.5
$ws.Range("I126").Value = 3485.6
$ws.Range("K126").Value = 10456.8
$ws.Range("M126").Value = -7986.799999999999
# row 132
$ws.Range("H132").Value = 4487.758
$ws.Range("I132").Value = 3426.3076
$ws.Range("J132").Value = 8430.286
$ws.Range("K132").Value = 10278.9228
$ws.Range("L132").Value = 25290.858
$ws.Range("M132").Value = -7748.9228
$ws.Range("N132").Value = -30350.858

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 5846.3076
$ws.Range("I22").Value = 925.125
$ws.Range("J22").Value = 13720.2
$ws.Range("K22").Value = 925.125
$ws.Range("L22").Value = 13720.2
$ws.Range("M22").Value = -630.125
$ws.Range("N22").Value = -14310.2
# row 27
$ws.Range("H27").Value = 5846.3076
$ws.Range("I27").Value = 925.125
$ws.Range("J27").Value = 13720.2
$ws.Range("K27").Value = 925.125
$ws.Range("L27").Value = 13720.2
$ws.Range("M27").Value = -818.125
$ws.Range("N27").Value = -13934.2
# row 40
$ws.Range("H40").Value = 5785.6875
$ws.Range("I40").Value = 5277.1377
$ws.Range("J40").Value = 10701.667
$ws.Range("K40").Value = 5277.1377
$ws.Range("L40").Value = 10701.667
$ws.Range("M40").Value = -5141.1377
$ws.Range("N40").Value = -10973.667
# row 68
$ws.Range("H68").Value = 3919
$ws.Range("I68").Value = 3362.6365
$ws.Range("K68").Value = 3362.6365
$ws.Range("M68").Value = -2613.6365
# row 71
$ws.Range("H71").Value = 3919
$ws.Range("I71").Value = 3362.6365
$ws.Range("K71").Value = 16813.1825
$ws.Range("M71").Value = -13069.1825

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 3306.077
$ws.Range("I81").Value = 1313.2858
$ws.Range("K81").Value = 2626.5716
$ws.Range("M81").Value = -1565.5716
# row 84
$ws.Range("H84").Value = 3306.077
$ws.Range("I84").Value = 1313.2858
$ws.Range("K84").Value = 13132.858
$ws.Range("M84").Value = -7828.858
# row 103
$ws.Range("H103").Value = 29633.334
$ws.Range("J103").Value = 29633.334
$ws.Range("L103").Value = 29633.334
$ws.Range("N103").Value = -31977.334
# row 107
$ws.Range("H107").Value = 1030.4
$ws.Range("I107").Value = 879.75
$ws.Range("K107").Value = 2639.25
$ws.Range("M107").Value = -719.25
